$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.286.41'
$ws.Cells.Item(2, 5).Value = '  +0.87%  '

$ws.Cells.Item(3, 4).Value = '1.657.17'
$ws.Cells.Item(3, 5).Value = '  +0.23%  '

$ws.Cells.Item(4, 5).Value = '  +0.61%  '

$ws.Cells.Item(5, 4).Value = '219.34'
$ws.Cells.Item(5, 5).Value = '  +2.29%  '

$ws.Cells.Item(6, 4).Value = '0.5221'
$ws.Cells.Item(6, 5).Value = '  -0.34%  '

$ws.Cells.Item(7, 4).Value = '1.007'
$ws.Cells.Item(7, 5).Value = '  +0.63%  '

$ws.Cells.Item(8, 4).Value = '0.2652'
$ws.Cells.Item(8, 5).Value = '  +1.11%  '

$ws.Cells.Item(9, 4).Value = '0.06324'
$ws.Cells.Item(9, 5).Value = '  -0.71%  '

$ws.Cells.Item(10, 4).Value = '21.39'
$ws.Cells.Item(10, 5).Value = '  +3.00%  '

$ws.Cells.Item(11, 4).Value = '0.07760'
$ws.Cells.Item(11, 5).Value = '  +0.24%  '

$ws.Cells.Item(12, 4).Value = '1.667.79'
$ws.Cells.Item(12, 5).Value = '  +1.01%  '

$ws.Cells.Item(13, 4).Value = '4.439'
$ws.Cells.Item(13, 5).Value = '  -0.18%  '

$ws.Cells.Item(14, 4).Value = '0.5485'
$ws.Cells.Item(14, 5).Value = '  -0.40%  '

$ws.Cells.Item(15, 4).Value = '0.0₅8232'
$ws.Cells.Item(15, 5).Value = '  -0.39%  '

$ws.Cells.Item(16, 4).Value = '65.04'
$ws.Cells.Item(16, 5).Value = '  +0.23%  '

$ws.Cells.Item(17, 4).Value = '26.326.77'
$ws.Cells.Item(17, 5).Value = '  +1.04%  '

$ws.Cells.Item(18, 5).Value = '  +0.45%  '

$ws.Cells.Item(19, 4).Value = '4.697'
$ws.Cells.Item(19, 5).Value = '  -0.99%  '

$ws.Cells.Item(20, 4).Value = '191.92'
$ws.Cells.Item(20, 5).Value = '  +0.85%  '

$ws.Cells.Item(21, 4).Value = '10.22'
$ws.Cells.Item(21, 5).Value = '  -0.37%  '

$ws.Cells.Item(22, 4).Value = '6.214'
$ws.Cells.Item(22, 5).Value = '  -2.15%  '

$ws.Cells.Item(23, 5).Value = '  +0.79%  '

$ws.Cells.Item(24, 4).Value = '138.88'
$ws.Cells.Item(24, 5).Value = '  -3.00%  '

$ws.Cells.Item(25, 4).Value = '0.1253'
$ws.Cells.Item(25, 5).Value = '  +0.45%  '

$ws.Cells.Item(26, 4).Value = '7.313'
$ws.Cells.Item(26, 5).Value = '  -1.17%  '

$ws.Cells.Item(27, 5).Value = '  +0.65%  '

$ws.Cells.Item(28, 4).Value = '1.424'
$ws.Cells.Item(28, 5).Value = '  +0.48%  '

$ws.Cells.Item(29, 4).Value = '0.06069'
$ws.Cells.Item(29, 5).Value = '  +2.15%  '

$ws.Cells.Item(30, 4).Value = '1.288'
$ws.Cells.Item(30, 5).Value = '  +2.16%  '

$ws.Cells.Item(31, 4).Value = '3.560'
$ws.Cells.Item(31, 5).Value = '  +3.54%  '

$ws.Cells.Item(32, 4).Value = '3.374'
$ws.Cells.Item(32, 5).Value = '  -1.18%  '

$ws.Cells.Item(33, 4).Value = '1.663'
$ws.Cells.Item(33, 5).Value = '  +0.70%  '

$ws.Cells.Item(34, 4).Value = '0.9902'
$ws.Cells.Item(34, 5).Value = '  -0.56%  '

$ws.Cells.Item(35, 5).Value = '  +1.18%  '

$ws.Cells.Item(36, 4).Value = '2.772'
$ws.Cells.Item(36, 5).Value = '  +0.54%  '

$ws.Cells.Item(37, 4).Value = '0.5980'
$ws.Cells.Item(37, 5).Value = '  +6.36%  '

$ws.Cells.Item(38, 4).Value = '0.01604'
$ws.Cells.Item(38, 5).Value = '  +0.03%  '

$ws.Cells.Item(39, 4).Value = '5.985'
$ws.Cells.Item(39, 5).Value = '  +1.91%  '

$ws.Cells.Item(40, 4).Value = '1.077.83'
$ws.Cells.Item(40, 5).Value = '  +4.88%  '

$ws.Cells.Item(41, 4).Value = '0.8564'
$ws.Cells.Item(41, 5).Value = '  -0.04%  '

$ws.Cells.Item(42, 5).Value = '  +0.37%  '

$ws.Cells.Item(43, 4).Value = '100.00'
$ws.Cells.Item(43, 5).Value = '  +0.81%  '

$ws.Cells.Item(44, 4).Value = '1.803.42'
$ws.Cells.Item(44, 5).Value = '  +0.29%  '

$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).Value = '57.45'
$ws.Cells.Item(45, 5).Value = '  +2.90%  '

$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).Value = '0.0₈107'
$ws.Cells.Item(46, 5).Value = '  -1.03%  '

$ws.Cells.Item(47, 2).Value = 'Frax'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(47, 4).Value = '1.003'
$ws.Cells.Item(47, 5).Value = '  +0.10%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).Value = '8.088'
$ws.Cells.Item(48, 5).Value = '  +0.52%  '

$ws.Cells.Item(49, 4).Value = '0.05196'
$ws.Cells.Item(49, 5).Value = '  +0.89%  '

$ws.Cells.Item(50, 4).Value = '1.471'
$ws.Cells.Item(50, 5).Value = '  +6.14%  '

$ws.Cells.Item(51, 4).Value = '0.4236'
$ws.Cells.Item(51, 5).Value = '  +0.63%  '

